$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Razem" (Total) header in E1, copying the bold/centered header style
# used by A1:C1, then add a SUM formula in E2 that totals column B.
$ws.Range("C1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Razem"
$ws.Range("E2").Formula = "=SUM(B:B)"

# Match the saved selection state from the diff.
[void]$ws.Range("E3").Select()
